$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): a new "Jan_2026" column is inserted before the
# existing "Dec_2025"/"Nov_2025" columns, which shift right by one column and the
# oldest month ("Oct_2025") is dropped. "ISIN"/"Stock Name"/"Mutual Fund"/"MoM"/"QoQ"
# stay where they are.
$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"

# --- Refreshed holdings data (rows 2-20), already sorted descending by the new
# Jan_2026 column, exactly as produced by the quant engine refresh.
$data = @(
    @('INE406A01037', 'Aurobindo Pharma Limited', 'quant Manufacturing Fund', 10.101106, 9.993342, 9.893630999999999, 0.1077639999999995, 0.2074750000000005),
    @('INE775A01035', 'Samvardhana Motherson International Ltd', 'quant Manufacturing Fund', 9.837403, 10.041152, 9.29668, -0.2037490000000002, 0.5407229999999998),
    @('INE423A01024', 'Adani Enterprises Limited', 'quant Manufacturing Fund', 9.614039, 9.733378, 0, -0.1193390000000001, 9.614039),
    @('INE768C01028', 'Zydus Wellness Ltd', 'quant Manufacturing Fund', 8.306566, 7.692854, 6.954599, 0.6137120000000005, 1.351967),
    @('INE045A01017', 'Ador Welding Limited', 'quant Manufacturing Fund', 6.21963, 5.795094, 5.630038, 0.4245360000000007, 0.5895920000000006),
    @('INE917I01010', 'Bajaj Auto Limited', 'quant Manufacturing Fund', 6.031482, 5.362379, 4.972071, 0.6691029999999998, 1.059411),
    @('INE364U01010', 'Adani Green Energy Limited', 'quant Manufacturing Fund', 5.856657, 6.370469, 7.284482, -0.5138119999999997, -1.427824999999999),
    @('INE206N01018', 'Ravindra Energy Limited', 'quant Manufacturing Fund', 5.600932, 5.541504, 5.019338, 0.05942800000000048, 0.5815939999999999),
    @('INE942C01045', 'Gujarat Themis Biosyn Ltd', 'quant Manufacturing Fund', 5.053164, 6.308265, 5.558495, -1.255101, -0.505331),
    @('INE180C01042', 'Capri Global Capital Limited', 'quant Manufacturing Fund', 4.176886, 3.968793, 3.899058, 0.2080929999999999, 0.2778279999999995),
    @('INE931S01010', 'Adani Energy Solutions Limited', 'quant Manufacturing Fund', 3.793646, 3.977906, 3.676665, -0.1842600000000001, 0.116981),
    @('INE331A01037', 'The Ramco Cements Limited', 'quant Manufacturing Fund', 3.159637, 0, 0, 3.159637, 3.159637),
    @('INE019C01026', 'Himadri Speciality Chemical Limited', 'quant Manufacturing Fund', 0, 0, 2.810666, 0, -2.810666),
    @('INE290A01027', 'Nahar Spinning Mills Limited', 'quant Manufacturing Fund', 0, 0.577651, 0.573519, -0.577651, -0.573519),
    @('INE669C01036', 'Tech Mahindra Limited', 'quant Manufacturing Fund', 0, 1.336715, 0, -1.336715, 0),
    @('INE769A01020', 'Aarti Industries Ltd', 'quant Manufacturing Fund', 0, 0, 2.815853, 0, -2.815853),
    @('INE171Z01026', 'Bharat Dynamics Limited', 'quant Manufacturing Fund', 0, 5.330176, 0, -5.330176, 0),
    @('INE0BS701011', 'Premier Energies Limited', 'quant Manufacturing Fund', 0, 0, 5.782484, 0, -5.782484),
    @('INE002A01018', 'Reliance Industries Limited', 'quant Manufacturing Fund', 0, 0, 1.718982, 0, -1.718982)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# --- The refresh only keeps 19 holdings now (rows 2-20); the previous rows
# 21-26 held stocks that dropped out of the portfolio entirely, so remove them.
$ws.Rows("21:26").Delete()
